$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 ---
$ws.Range("A6").Value = "'DW24LQH232"
$ws.Range("B6").Value = 12
$ws.Range("C6").Value = "'7.4.2015 г. 00:00:00 ч."
$ws.Range("D6").Value = "'IrregularExpense"
$ws.Range("E6").Value = "'test excel"

# --- Row 7 ---
$ws.Range("A7").Value = "'7OFMUXDQU9"
$ws.Range("B7").Value = 34
$ws.Range("C7").Value = "'7.4.2015 г. 00:00:00 ч."
$ws.Range("D7").Value = "'RegularExpense"
$ws.Range("E7").Value = "'2fsd"

# Give the numeric Amount cells (column B) the same cell style (quote-prefix
# formatting) as the rest of the data rows, like the existing B2:B5 cells,
# without disturbing their numeric value/type.
$ws.Range("B5").Copy()
$ws.Range("B6").PasteSpecial(-4122)
$ws.Range("B7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
